# Updates cryptos list values (Price column D and Volume(1h) column E)
# plus a row swap between Polkadot (row 12->13) and WrappedEther (row 13->12),
# matching the "Updated cryptos list" GitHub Actions commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.543.18"
$ws.Range("E2").Value = "  +0.00%  "
$ws.Range("D3").Value = "1.922.43"
$ws.Range("E3").Value = "  +0.47%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.013"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.44%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "326.05"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.011"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.33%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4821"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4070"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08241"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.88%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.010"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "23.57"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.83%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.085"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.35%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.886.42"
$ws.Range("E13").Value = "  -1.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.275"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +2.31%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.70"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.42%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06862"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.77%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.013"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.32%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001039"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.23%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.61"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.65%  "
$ws.Range("E20").Value = "  +0.30%  "
$ws.Range("D21").Value = "29.556.73"
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.688"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.90"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.51%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.184"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.26%  "
$ws.Range("D25").Value = "2.114.49"
$ws.Range("E25").Value = "  -1.17%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "155.88"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.23%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.476"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.92%  "
$ws.Range("E28").Value = "  -0.21%  "
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "120.68"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.98%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.017"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.92%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09650"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.83%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.634"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.85%  "
$ws.Range("E34").Value = "  -0.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.375"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06366"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +4.12%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02298"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.31%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.190"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.23%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5954"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.34%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.78"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.896"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.32%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1851"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.477"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.25%  "
$ws.Range("E44").Value = "  -0.36%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.47"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.34%  "
$ws.Range("E46").Value = "  -2.90%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5573"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.22%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.951"
$ws.Range("D48").ClearFormats()
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "118.93"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.95%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.435"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +3.36%  "
$ws.Range("E51").Value = "  -0.57%  "